$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EmployeeData")

# Row 2 - Tom Holland
$ws.Range("A2").Value = "Tom"
$ws.Range("B2").Value = "Sr"
$ws.Range("C2").Value = "Holland"
$ws.Range("D2").Value = "C:\Users\quytr\Desktop\Batch12\imagebatch12.jpg"
$ws.Range("E2").Value = "tom123ho"
$ws.Range("F2").Value = "Hum@nhrm123"

# Row 3 - Jackie Chan
$ws.Range("A3").Value = "Jackie"
$ws.Range("B3").Value = "Jr"
$ws.Range("C3").Value = "Chan"
$ws.Range("D3").Value = "C:\Users\quytr\Desktop\Batch12\imagebatch12.jpg"
$ws.Range("E3").Value = "jack8943ch"
$ws.Range("F3").Value = "Hum@nhrm123"

# Row 4 - Star Lord
$ws.Range("A4").Value = "Star"
$ws.Range("B4").Value = "Cool"
$ws.Range("C4").Value = "Lord"
$ws.Range("D4").Value = "C:\Users\quytr\Desktop\Batch12\imagebatch12.jpg"
$ws.Range("E4").Value = "star130lord"
$ws.Range("F4").Value = "Hum@nhrm123"

# Update the active selection to B2
$ws.Range("B2").Select()
